$d = $word.ActiveDocument

# Locate the paragraph that holds the "LOQ4088: ..." requirement line and the
# paragraph holding the "© 2020 . Contact: ..." footer line, then remove the
# paragraphs in between (a blank line, the "Ver no Jupiter ..." line, and the
# "© 2020 ..." line itself) together with their paragraph marks, while leaving
# the single blank paragraph that precedes the final page-break paragraph
# untouched.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    if ($text -like "LOQ4088:*") {
        $startPara = $para
    }
    elseif ($text -like "*Contact: luizeleno@usp.br*") {
        $endPara = $para
        break
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $r = $d.Range($startPara.Range.End, $endPara.Range.End)
    $r.Delete()
}
